# Update the "Server side and API" / "Client side and Security" bullet
# lines on the second slide (placeholder content box) to reflect the new
# wording from the commit ("Server side, API and Security" /
# "Client side and Web Site").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

function Set-SubstringText($Shape, $OldText, $NewText) {
    $full = $Shape.TextFrame.TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Could not find text '$OldText' in shape '$($Shape.Name)'"
    }

    # PowerPoint TextRange.Characters is 1-based.
    $run = $Shape.TextFrame.TextRange.Characters($idx + 1, $OldText.Length)
    $run.Text = $NewText
}

Set-SubstringText $sh "Server side and API" "Server side, API and Security"
Set-SubstringText $sh "Client side and Security" "Client side and Web Site"
